$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 1327

$data = @(
    @(45156, 26040.81931205392, 4, 0, 0, 0.9186590956221448, 5.08),
    @(45157, 26115.91295705558, 1, 0, 1, 0.9208352491360976, 5.08),
    @(45158, 26137.35350149368, 2, 0, 0, 0.9267152860013572, 5.08),
    @(45159, 26086.51127135073, 4, 0, 0, 0.9173868726084622, 5.08),
    @(45160, 25826.03364121723, 2, 0, 1, 0.9216308215777896, 5.08),
    @(45161, 26641.12006805852, 3, 0, 0, 0.916416413647274, 5.08),
    @(45162, 26021.08163211241, 3, 0, 0, 0.924915071315602, 5.08),
    @(45163, 25949.31275796248, 1, 0, 0, 0.9253015746700952, 5.08),
    @(45164, 26025.2295197059, 1, 0, 0, 0.9260709106473284, 5.08),
    @(45165, 26112.34304465044, 2, 0, 0, 0.9291095803036642, 5.08),
    @(45166, 25975.58230041972, 2, 0, 0, 0.9245288456351862, 5.08),
    @(45167, 27852.18330418598, 5, 0, 0, 0.9182388455392472, 5.08),
    @(45168, 27251.96914888735, 4, 0, 1, 0.9140910262678168, 5.08),
    @(45169, 26179.07181383978, 4, 0, 1, 0.9183313397383228, 5.08),
    @(45170, 25581.65193106596, 3, 0, 0, 0.9275432091775184, 5.08),
    @(45171, 25845.21631827925, 3, 0, 0, 0.922686074170466, 5.08),
    @(45172, 25991.65745373537, 3, 0, 0, 0.9284101405876204, 5.08),
    @(45173, 25871.5364941652, 2, 0, 1, 0.9265092605722368, 5.08),
    @(45174, 25706.43383971701, 2, 0, 0, 0.9320291692234548, 5.08),
    @(45175, 25643.30464824743, 9, 0, 0, 0.9321445941858784, 5.08),
    @(45176, 25871.77707088979, 4, 0, 1, 0.9334941303910536, 5.08),
    @(45177, 25880.11091714333, 2, 0, 0, 0.9327919389566486, 5.08),
    @(45178, 25857.72787564452, 3, 0, 0, 0.9327130559698408, 5.08)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("A$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}
